$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.016169548034668
$ws.Range("B1").Value = 2.035089731216431
$ws.Range("C1").Value = 3.685343265533447
$ws.Range("D1").Value = 1.910654187202454
$ws.Range("E1").Value = 0.398711770772934
